$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.475.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.21%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.911.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.48%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'239.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.47%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.9992"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.11%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4772"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.74%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.2847"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.91%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.06706"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.74%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'19.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.42%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'103.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.01%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07754"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'1.917.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.28%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.189"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -3.11%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.6686"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.89%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'276.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.35%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'30.497.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.20%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("E18").Value = "'  -0.13%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.000007488"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.05%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'12.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.73%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'5.381"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.49%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").Value = "'BinanceUSD"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.9992"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.11%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").Value = "'Chainlink"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'6.293"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.91%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = "'Cosmos"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'9.341"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.49%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").Value = "'Monero"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'166.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.25%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").Value = "'EthereumClassic"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'19.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.92%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").Value = "'LidoDAOToken"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = "'2.076"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -3.61%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = "'Toncoin"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'1.380"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.74%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").Value = "'Stellar"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'0.09982"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.74%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = "'Filecoin"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'4.601"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.58%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "'PancakeSwap"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'1.510"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.03%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'InternetComputer(DFINITY)"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'4.259"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.66%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = "'Hedera"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'0.04695"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.74%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = "'ImmutableX"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.7284"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.55%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").Value = "'ARBITRUM"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'1.113"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -3.39%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").Value = "'HuobiToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'2.708"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.11%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'VeChain"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.01906"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.52%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'MXToken"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'2.608"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.84%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "'FraxShare"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'6.378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.22%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "'Aave"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'74.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.75%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "'RenderToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'1.959"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -6.18%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "'TrustWalletToken"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.8613"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.76%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "'Quant"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'106.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.61%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "'TheSandbox"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.4266"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.66%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'PaxDollar"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'0.9985"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.04%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'Aptos"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'7.429"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.94%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'Maker"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'950.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.28%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").Value = "'Algorand"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.1210"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.55%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").Value = "'Elrond"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'34.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.10%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").Value = "'Cronos"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.05796"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.40%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'8.738"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.30%  "
$ws.Range("E51").Style = "Normal"
